# Testing for RFC implemented
# Adds a new worksheet "CompleteRFC5" at the end of the workbook, as a
# duplicate of the "CompleteRFC4" sheet (same P1-P5/Promedio/Accuracy
# layout, values and formatting), per the RFC testing results table.

$wb = $excel.ActiveWorkbook

# Remember the sheet that is active before we start, so the workbook's
# active-tab selection is left untouched by the new sheet creation.
$origActive = $wb.ActiveSheet

# Source sheet to clone: the most recently added "CompleteRFC*" results sheet.
$sourceSheet = $wb.Worksheets.Item("CompleteRFC4")

# Insert the clone right after the last worksheet in the workbook so the
# new tab lands at the very end (position 10).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$sourceSheet.Copy([System.Reflection.Missing]::Value, $lastSheet)

$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "CompleteRFC5"

# Restore the originally active sheet/tab selection.
$origActive.Activate()

Write-Host "Added worksheet '$($newSheet.Name)' (now $($wb.Worksheets.Count) sheets total)."
